# Apply the "E suite" workbook edit: append TestCase_E42 as row 43 on the
# "Test Cases" sheet, then move the active selection to A2 (matching the
# saved sheetView/selection seen in the target workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# --- Copy the formatting from the most similar existing rows -------------
# Row 40 already has the "TCID / Jira(wrap) / Description(wrap) / Y / PASS"
# shape with a row height auto-expanded for wrapped text; row 30's
# Description cell (column C) carries the alternate shaded wrap style that
# the new row needs, so borrow it for C43.
$ws.Range("A40:E40").Copy()
$ws.Range("A43:E43").PasteSpecial(-4122)
$ws.Range("C30").Copy()
$ws.Range("C43").PasteSpecial(-4122)

# --- Cell values for the new test case ------------------------------------
# (Description is assigned before Jira id so the new shared-string table
# entries land in the same order as the saved workbook: TCID, Description,
# then Jira id.)
$ws.Range("A43").Value = "TestCase_E42"
$ws.Range("C43").Value = "Verify that user is able to watch an article to a particular watchlist from notification in home page||Verify that user is able to unwatch an article from watchlist from notification in home page"
$ws.Range("B43").Value = "OPQA-298`n||OPQA-304"
$ws.Range("D43").Value = "Y"
$ws.Range("E43").Value = "PASS"

# Match the row height used by the other multi-line rows in the sheet.
$ws.Rows.Item(43).RowHeight = 30

$excel.CutCopyMode = 0

# --- View state -------------------------------------------------------
# The saved file scrolls back to the top and selects A2 instead of the
# previous C40 / topLeftCell A31 scroll position.
$ws.Activate()
$ws.Range("A2").Select()
